$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3181.1875
$ws.Range("I38").Value = 3181.1875
$ws.Range("K38").Value = 9543.5625
$ws.Range("M38").Value = -9171.5625
$ws.Range("H40").Value = 4217.2173
$ws.Range("J40").Value = 5351.909
$ws.Range("L40").Value = 5351.909
$ws.Range("N40").Value = -5701.909
$ws.Range("H111").Value = 3210.9333
$ws.Range("I111").Value = 2790.4614
$ws.Range("J111").Value = 5944
$ws.Range("K111").Value = 8371.3842
$ws.Range("L111").Value = 17832
$ws.Range("M111").Value = -5304.3842
$ws.Range("N111").Value = -23966
$ws.Range("H137").Value = 3878.8367
$ws.Range("I137").Value = 4087.6924
$ws.Range("K137").Value = 12263.0772
$ws.Range("M137").Value = -9713.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("H32").Value = 11498447
$ws.Range("I32").Value = 12824941
$ws.Range("K32").Value = 12824941
$ws.Range("M32").Value = -12824654
$ws.Range("H45").Value = 3301.8667
$ws.Range("I45").Value = 2972.8
$ws.Range("K45").Value = 2972.8
$ws.Range("M45").Value = -2595.8
$ws.Range("H61").Value = 1580.6451
$ws.Range("I61").Value = 1566.7916
$ws.Range("J61").Value = 1628.1428
$ws.Range("K61").Value = 1566.7916
$ws.Range("L61").Value = 1628.1428
$ws.Range("M61").Value = -1354.7916
$ws.Range("N61").Value = -2052.1428
$ws.Range("H74").Value = 1597.08
$ws.Range("I74").Value = 1361
$ws.Range("J74").Value = 3328.3333
$ws.Range("K74").Value = 1361
$ws.Range("L74").Value = 3328.3333
$ws.Range("M74").Value = -487
$ws.Range("N74").Value = -5076.3333
$ws.Range("H77").Value = 1597.08
$ws.Range("I77").Value = 1361
$ws.Range("J77").Value = 3328.3333
$ws.Range("K77").Value = 6805
$ws.Range("L77").Value = 16641.6665
$ws.Range("M77").Value = -2437
$ws.Range("N77").Value = -25377.6665
$ws.Range("H102").Value = 27781252
$ws.Range("I102").Value = 3776.7778
$ws.Range("K102").Value = 3776.7778
$ws.Range("M102").Value = -2154.7778
$ws.Range("H132").Value = 1645.3871
$ws.Range("I132").Value = 1083.3636
$ws.Range("J132").Value = 3019.2222
$ws.Range("K132").Value = 3250.0908
$ws.Range("L132").Value = 9057.6666
$ws.Range("M132").Value = -720.0907999999999
$ws.Range("N132").Value = -14117.6666
$ws.Range("H136").Value = 1580.6451
$ws.Range("I136").Value = 1566.7916
$ws.Range("J136").Value = 1628.1428
$ws.Range("K136").Value = 4700.3748
$ws.Range("L136").Value = 4884.428400000001
$ws.Range("M136").Value = -2150.3748
$ws.Range("N136").Value = -9984.428400000001
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 69003.87
$ws.Range("I20").Value = 1388.375
$ws.Range("J20").Value = 146278.72
$ws.Range("K20").Value = 1388.375
$ws.Range("L20").Value = 146278.72
$ws.Range("M20").Value = -1141.375
$ws.Range("N20").Value = -146772.72
$ws.Range("H86").Value = 2856.9778
$ws.Range("I86").Value = 2309.7646
$ws.Range("J86").Value = 3189.2144
$ws.Range("K86").Value = 2309.7646
$ws.Range("L86").Value = 3189.2144
$ws.Range("M86").Value = -1186.7646
$ws.Range("N86").Value = -5435.2144
$ws.Range("H89").Value = 2856.9778
$ws.Range("I89").Value = 2309.7646
$ws.Range("J89").Value = 3189.2144
$ws.Range("K89").Value = 11548.823
$ws.Range("L89").Value = 15946.072
$ws.Range("M89").Value = -5932.823
$ws.Range("N89").Value = -27178.072
$ws.Range("H94").Value = 5001488
$ws.Range("I94").Value = 1756.1666
$ws.Range("K94").Value = 1756.1666
$ws.Range("M94").Value = -1305.1666
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H134").Value = 2588.5615
$ws.Range("I134").Value = 1811.5883
$ws.Range("K134").Value = 5434.7649
$ws.Range("M134").Value = -2899.7649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 27398.6
$ws.Range("J28").Value = 31748.25
$ws.Range("L28").Value = 31748.25
$ws.Range("N28").Value = -32238.25
$ws.Range("H31").Value = 2177.6904
$ws.Range("I31").Value = 1714.091
$ws.Range("J31").Value = 2687.65
$ws.Range("K31").Value = 1714.091
$ws.Range("L31").Value = 2687.65
$ws.Range("M31").Value = -1419.091
$ws.Range("N31").Value = -3277.65
$ws.Range("H34").Value = 2177.6904
$ws.Range("I34").Value = 1714.091
$ws.Range("J34").Value = 2687.65
$ws.Range("K34").Value = 1714.091
$ws.Range("L34").Value = 2687.65
$ws.Range("M34").Value = -1512.091
$ws.Range("N34").Value = -3091.65
$ws.Range("H53").Value = 42622.5
$ws.Range("J53").Value = 42622.5
$ws.Range("L53").Value = 42622.5
$ws.Range("N53").Value = -43836.5
$ws.Range("H86").Value = 88393.39999999999
$ws.Range("I86").Value = 86249.5
$ws.Range("K86").Value = 86249.5
$ws.Range("M86").Value = -85126.5
$ws.Range("H89").Value = 88393.39999999999
$ws.Range("I89").Value = 86249.5
$ws.Range("K89").Value = 431247.5
$ws.Range("M89").Value = -425631.5
$ws.Range("H92").Value = 79997.5
$ws.Range("J92").Value = 79997.5
$ws.Range("L92").Value = 79997.5
$ws.Range("N92").Value = -84989.5
$ws.Range("H93").Value = 30549.777
$ws.Range("I93").Value = 12491.5
$ws.Range("J93").Value = 66666.336
$ws.Range("K93").Value = 12491.5
$ws.Range("L93").Value = 66666.336
$ws.Range("M93").Value = -10619.5
$ws.Range("N93").Value = -70410.336
$ws.Range("H95").Value = 3178.1667
$ws.Range("J95").Value = 3178.1667
$ws.Range("L95").Value = 3178.1667
$ws.Range("N95").Value = -8670.1667
$ws.Range("H96").Value = 19244.8
$ws.Range("J96").Value = 19244.8
$ws.Range("L96").Value = 19244.8
$ws.Range("N96").Value = -24736.8
$ws.Range("H97").Value = 49999.668
$ws.Range("J97").Value = 49999.668
$ws.Range("L97").Value = 49999.668
$ws.Range("N97").Value = -51981.668
$ws.Range("H107").Value = 8075.3335
$ws.Range("I107").Value = 435.10526
$ws.Range("J107").Value = 26220.875
$ws.Range("K107").Value = 435.10526
$ws.Range("L107").Value = 26220.875
$ws.Range("M107").Value = 1484.89474
$ws.Range("N107").Value = -30060.875
$ws.Range("H134").Value = 2199.3547
$ws.Range("I134").Value = 2038.1852
$ws.Range("K134").Value = 6114.5556
$ws.Range("M134").Value = -3579.5556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 31469274
$ws.Range("J4").Value = 670.5
$ws.Range("L4").Value = 2011.5
$ws.Range("N4").Value = -2235.5
$ws.Range("H56").Value = 5797.2
$ws.Range("I56").Value = 5797.2
$ws.Range("K56").Value = 5797.2
$ws.Range("M56").Value = -5267.2
$ws.Range("H132").Value = 1980.5
$ws.Range("J132").Value = 1494.6666
$ws.Range("L132").Value = 13451.9994
$ws.Range("N132").Value = -18511.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6616.077
$ws.Range("I113").Value = 3882
$ws.Range("J113").Value = 8324.875
$ws.Range("K113").Value = 3882
$ws.Range("L113").Value = 8324.875
$ws.Range("M113").Value = -1712
$ws.Range("N113").Value = -12664.875
$ws.Range("H122").Value = 4557.2593
$ws.Range("I122").Value = 3427.111
$ws.Range("K122").Value = 10281.333
$ws.Range("M122").Value = -7831.332999999999
$ws.Range("H132").Value = 2360.8281
$ws.Range("I132").Value = 2082.1428
$ws.Range("J132").Value = 3271.2
$ws.Range("K132").Value = 6246.428400000001
$ws.Range("L132").Value = 9813.599999999999
$ws.Range("M132").Value = -3716.428400000001
$ws.Range("N132").Value = -14873.6
$ws.Range("H134").Value = 72016.78
$ws.Range("J134").Value = 72016.78
$ws.Range("L134").Value = 216050.34
$ws.Range("N134").Value = -221120.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13360.1
$ws.Range("I40").Value = 16384.285
$ws.Range("J40").Value = 6303.6665
$ws.Range("K40").Value = 16384.285
$ws.Range("L40").Value = 6303.6665
$ws.Range("M40").Value = -16248.285
$ws.Range("N40").Value = -6575.6665
$ws.Range("H68").Value = 3922.7222
$ws.Range("I68").Value = 3890.2222
$ws.Range("K68").Value = 3890.2222
$ws.Range("M68").Value = -3141.2222
$ws.Range("H71").Value = 3922.7222
$ws.Range("I71").Value = 3890.2222
$ws.Range("K71").Value = 19451.111
$ws.Range("M71").Value = -15707.111
$ws.Range("H80").Value = 101857.14
$ws.Range("I80").Value = 95000
$ws.Range("J80").Value = 103000
$ws.Range("K80").Value = 95000
$ws.Range("L80").Value = 103000
$ws.Range("M80").Value = -93877
$ws.Range("N80").Value = -105246
$ws.Range("H83").Value = 101857.14
$ws.Range("I83").Value = 95000
$ws.Range("J83").Value = 103000
$ws.Range("K83").Value = 285000
$ws.Range("L83").Value = 309000
$ws.Range("M83").Value = -279384
$ws.Range("N83").Value = -320232
$ws.Range("H132").Value = 3116.7285
$ws.Range("I132").Value = 2056.1667
$ws.Range("K132").Value = 6168.500100000001
$ws.Range("M132").Value = -3638.500100000001
$ws.Range("H136").Value = 4244.423
$ws.Range("I136").Value = 4244.423
$ws.Range("K136").Value = 12733.269
$ws.Range("M136").Value = -10183.269

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3012.4285
$ws.Range("I122").Value = 2483.28
$ws.Range("J122").Value = 4335.3
$ws.Range("K122").Value = 7449.84
$ws.Range("L122").Value = 13005.9
$ws.Range("M122").Value = -4999.84
$ws.Range("N122").Value = -17905.9
$ws.Range("H126").Value = 1362.2858
$ws.Range("I126").Value = 950.8
$ws.Range("K126").Value = 2852.4
$ws.Range("M126").Value = -382.3999999999996
$ws.Range("H132").Value = 1581.878
$ws.Range("I132").Value = 1207.7059
$ws.Range("K132").Value = 3623.1177
$ws.Range("M132").Value = -1093.1177
